{"js": "// Update every cell in the first (and only) table with the new values,\n// preserving table structure, cell formatting, and the leading date paragraph.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, row-major, 20 rows x 5 columns (matches the diff order exactly).\nconst newValues = [\n  [\"24+55=\", \"45+41=\", \"40+16=\", \"22-1=\", \"78-51=\"],\n  [\"22+73=\", \"93-1=\", \"89-87=\", \"27-19=\", \"1+63=\"],\n  [\"12-7=\", \"79-16=\", \"76-60=\", \"15+34=\", \"96-15=\"],\n  [\"93-49=\", \"40+44=\", \"22+32=\", \"1+55=\", \"86-33=\"],\n  [\"86+6=\", \"97-7=\", \"12+71=\", \"50-35=\", \"84-14=\"],\n  [\"22+17=\", \"78-16=\", \"31+3=\", \"49-46=\", \"60-10=\"],\n  [\"81-37=\", \"59+2=\", \"27-1=\", \"43+5=\", \"24+19=\"],\n  [\"17+82=\", \"65+4=\", \"23-16=\", \"74-29=\", \"38-5=\"],\n  [\"45-30=\", \"7+79=\", \"11+18=\", \"9+29=\", \"8+64=\"],\n  [\"99-40=\", \"53+35=\", \"4+5=\", \"93-79=\", \"20-9=\"],\n  [\"10+82=\", \"56-26=\", \"84-35=\", \"1+67=\", \"1+9=\"],\n  [\"59+0=\", \"13+31=\", \"17+70=\", \"8+12=\", \"52+42=\"],\n  [\"77-34=\", \"91-6=\", \"49+36=\", \"7+30=\", \"76-11=\"],\n  [\"46-11=\", \"61-22=\", \"58+0=\", \"17+22=\", \"6+51=\"],\n  [\"99-42=\", \"37+53=\", \"72-68=\", \"48+21=\", \"88-86=\"],\n  [\"5+52=\", \"95-80=\", \"23+0=\", \"61-28=\", \"15-14=\"],\n  [\"36+2=\", \"83+14=\", \"54+38=\", \"56-56=\", \"18+51=\"],\n  [\"68-43=\", \"86-57=\", \"44+43=\", \"23+53=\", \"19-18=\"],\n  [\"12+50=\", \"30+55=\", \"5+17=\", \"19+14=\", \"33-2=\"],\n  [\"5+51=\", \"12+4=\", \"70+5=\", \"46+40=\", \"35+37=\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update every cell in the first (and only) table with the new values,\n# preserving table structure, cell formatting, and the leading date paragraph.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New values, row-major, 20 rows x 5 columns (matches the diff order exactly).\n$newValues = @(\n    @(\"24+55=\", \"45+41=\", \"40+16=\", \"22-1=\", \"78-51=\"),\n    @(\"22+73=\", \"93-1=\", \"89-87=\", \"27-19=\", \"1+63=\"),\n    @(\"12-7=\", \"79-16=\", \"76-60=\", \"15+34=\", \"96-15=\"),\n    @(\"93-49=\", \"40+44=\", \"22+32=\", \"1+55=\", \"86-33=\"),\n    @(\"86+6=\", \"97-7=\", \"12+71=\", \"50-35=\", \"84-14=\"),\n    @(\"22+17=\", \"78-16=\", \"31+3=\", \"49-46=\", \"60-10=\"),\n    @(\"81-37=\", \"59+2=\", \"27-1=\", \"43+5=\", \"24+19=\"),\n    @(\"17+82=\", \"65+4=\", \"23-16=\", \"74-29=\", \"38-5=\"),\n    @(\"45-30=\", \"7+79=\", \"11+18=\", \"9+29=\", \"8+64=\"),\n    @(\"99-40=\", \"53+35=\", \"4+5=\", \"93-79=\", \"20-9=\"),\n    @(\"10+82=\", \"56-26=\", \"84-35=\", \"1+67=\", \"1+9=\"),\n    @(\"59+0=\", \"13+31=\", \"17+70=\", \"8+12=\", \"52+42=\"),\n    @(\"77-34=\", \"91-6=\", \"49+36=\", \"7+30=\", \"76-11=\"),\n    @(\"46-11=\", \"61-22=\", \"58+0=\", \"17+22=\", \"6+51=\"),\n    @(\"99-42=\", \"37+53=\", \"72-68=\", \"48+21=\", \"88-86=\"),\n    @(\"5+52=\", \"95-80=\", \"23+0=\", \"61-28=\", \"15-14=\"),\n    @(\"36+2=\", \"83+14=\", \"54+38=\", \"56-56=\", \"18+51=\"),\n    @(\"68-43=\", \"86-57=\", \"44+43=\", \"23+53=\", \"19-18=\"),\n    @(\"12+50=\", \"30+55=\", \"5+17=\", \"19+14=\", \"33-2=\"),\n    @(\"5+51=\", \"12+4=\", \"70+5=\", \"46+40=\", \"35+37=\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
